# Update countries & provincias Spain
#
# 1) Move "Austria" up in the country ranking: row 38 becomes Austria (new
#    figures), row 39 becomes Japon (previous row-38 figures) and row 40
#    becomes Indonesia (previous row-39 figures).
# 2) Refresh the "last updated" timestamp string.
# 3) Update several countries' statistics (Casos totales, Nuevos casos,
#    Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Mayo de 2020 a las 12:35"

# --- Re-rank Austria / Japon / Indonesia, rows 38-40 ---
# Row 38: Austria (brand new figures)
$ws.Range("A38").Value = "Austria"
$ws.Range("B38").Value = 16058
$ws.Range("C38").Value = 61
$ws.Range("D38").Value = 14405
$ws.Range("E38").Value = 1027
$ws.Range("F38").Value = 54
$ws.Range("G38").Value = 2
$ws.Range("H38").Value = 626

# Row 39: Japon (takes the old row-38 figures)
$ws.Range("A39").Value = "Japon"
$ws.Range("B39").Value = 16049
$ws.Range("C39").Value = 0
$ws.Range("D39").Value = 8920
$ws.Range("E39").Value = 6451
$ws.Range("F39").Value = 243
$ws.Range("G39").Value = 0
$ws.Range("H39").Value = 678

# Row 40: Indonesia (takes the old row-39 figures)
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 16006
$ws.Range("C40").Value = 568
$ws.Range("D40").Value = 3518
$ws.Range("E40").Value = 11445
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 15
$ws.Range("H40").Value = 1043

# --- Row 25 (country rank 29) ---
$ws.Range("B25").Value = 30463
$ws.Range("C25").Value = 50
$ws.Range("E25").Value = 1493

# --- Row 53 (country rank 57) ---
$ws.Range("D53").Value = 6301
$ws.Range("E53").Value = 590

# --- Row 56 (country rank 60) ---
$ws.Range("B56").Value = 6593
$ws.Range("C56").Value = 81
$ws.Range("D56").Value = 3222
$ws.Range("E56").Value = 3182
$ws.Range("G56").Value = 1
$ws.Range("H56").Value = 189

# --- Row 99 (country rank 103) ---
$ws.Range("B99").Value = 1052
$ws.Range("C99").Value = 1
$ws.Range("D99").Value = 1009
